$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I14").Value = "sv"
$ws.Range("J14").Value = "Statement-opinion"

# Row 23: aa/Agree-Accept -> %/Uninterpretable
$ws.Range("I23").Value = "%"
$ws.Range("J23").Value = "Uninterpretable"

# Row 25: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I25").Value = "sv"
$ws.Range("J25").Value = "Statement-opinion"

# Row 37: b/Acknowledge (Backchannel) -> sd/Statement-non-opinion
$ws.Range("I37").Value = "sd"
$ws.Range("J37").Value = "Statement-non-opinion"

# Row 39: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I39").Value = "sv"
$ws.Range("J39").Value = "Statement-opinion"

# Row 48: qy/Yes-No-Question -> aa/Agree-Accept
$ws.Range("I48").Value = "aa"
$ws.Range("J48").Value = "Agree/Accept"

# Row 52: ba/Appreciation -> sd/Statement-non-opinion
$ws.Range("I52").Value = "sd"
$ws.Range("J52").Value = "Statement-non-opinion"

$wb.Save()
